$wb = $excel.ActiveWorkbook

# --- Rename "Flow Chart template tables" to "Flow Chart Template Tables" ---
$flowChartTablesSheet = $wb.Worksheets.Item("Flow Chart template tables")
$flowChartTablesSheet.Name = "Flow Chart Template Tables"

# --- Reorder merge cells on "Proposal (By Plan)" sheet (sheet4) ---
# Move G7:L7, N7:T7, G13:L13, N13:T13, C35:T35 to the front of the merge list.
$proposalByPlan = $wb.Worksheets.Item("Proposal (By Plan)")
$proposalByPlanOrder = @(
    "G7:L7", "N7:T7", "G13:L13", "N13:T13", "C35:T35",
    "C37:T37", "C39:T39", "C41:T41",
    "G19:L19", "N19:T19", "G25:L25", "N25:T25",
    "C31:T31", "C33:T33"
)
foreach ($r in $proposalByPlanOrder) {
    $proposalByPlan.Range($r).UnMerge()
}
foreach ($r in $proposalByPlanOrder) {
    $proposalByPlan.Range($r).Merge()
}

# --- Reorder merge cells on "Proposal" sheet (sheet6) ---
# Move C5:D5, H5:I5, J5:K5, L5:M5, F14:T14 to the end of the merge list.
$proposal = $wb.Worksheets.Item("Proposal")
$proposalOrder = @(
    "F16:T16", "F18:T18", "F20:T20",
    "O5:P5", "Q5:S5", "E5:G5",
    "H7:M7", "N7:T7", "F12:T12",
    "C5:D5", "H5:I5", "J5:K5", "L5:M5", "F14:T14"
)
foreach ($r in $proposalOrder) {
    $proposal.Range($r).UnMerge()
}
foreach ($r in $proposalOrder) {
    $proposal.Range($r).Merge()
}

# --- Reorder merge cells on "Flow Chart Template Tables" sheet (sheet7) ---
# Move C56:F56, G56:K56, L56:P56, C38:G38, H38:L38, M38:P38, C47:G47, H47:K47, L47:P47 to the front.
$flowChartTablesOrder = @(
    "C56:F56", "G56:K56", "L56:P56",
    "C38:G38", "H38:L38", "M38:P38",
    "C47:G47", "H47:K47", "L47:P47",
    "C2:F2", "G2:J2",
    "C29:F29", "G29:J29", "K29:O29",
    "K2:N2",
    "C11:G11", "H11:K11", "L11:O11",
    "C20:F20", "G20:K20", "L20:O20"
)
foreach ($r in $flowChartTablesOrder) {
    $flowChartTablesSheet.Range($r).UnMerge()
}
foreach ($r in $flowChartTablesOrder) {
    $flowChartTablesSheet.Range($r).Merge()
}

# --- Activate "Flow Chart Template Tables" sheet (moves tabSelected + activeTab) ---
$flowChartTablesSheet.Activate()
